$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the "amount of opened texture" values in columns N, O, P (rows 2-41)
# by 29 for each cell (fixes the hasValues issue referenced in the commit msg).
$values = @(
    @(36,11,42),
    @(47,40,1),
    @(17,52,59),
    @(23,35,47),
    @(1,32,39),
    @(1,26,24),
    @(50,7,32),
    @(55,13,8),
    @(12,11,58),
    @(10,44,15),
    @(58,23,42),
    @(12,28,4),
    @(58,21,46),
    @(35,31,49),
    @(50,44,13),
    @(12,40,32),
    @(13,53,19),
    @(10,16,37),
    @(9,21,30),
    @(5,53,16),
    @(20,38,29),
    @(2,15,12),
    @(31,37,51),
    @(25,29,59),
    @(38,54,29),
    @(3,59,7),
    @(7,15,20),
    @(3,28,51),
    @(29,11,31),
    @(60,28,44),
    @(50,32,38),
    @(41,7,26),
    @(45,53,29),
    @(13,46,45),
    @(50,15,2),
    @(55,10,49),
    @(39,38,34),
    @(34,54,59),
    @(22,50,13),
    @(9,36,56)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $triple = $values[$i]
    $ws.Cells.Item($row, 14).Value = $triple[0]
    $ws.Cells.Item($row, 15).Value = $triple[1]
    $ws.Cells.Item($row, 16).Value = $triple[2]
}

# Widen columns N:P (no longer auto-fit / bestFit)
$ws.Range("N1:P1").ColumnWidth = 9.75

# Restore the selection to N8
$ws.Range("N8").Select()
